# Enable AgGrid column grouping, and add import/export in customToolbar
# (data fixture update: fix Howard's birthday, add a new customer row "Billy")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1 ("customers") edits -------------------------------------------
# (sheet2 "UnusedSheet" keeps its two text cells untouched; switching the
#  active sheet below is the only change that affects it)

# Howard's birthday (E5) was wrong - fix the text value
$ws1.Range("E5").Value = "21/05/2002"

# Add a new customer row: Billy, 23, canDrinkAlcohol=TRUE, Beer, DATE(1940,4,28), 1.25
$ws1.Range("A6").Value = "Billy"
$ws1.Range("B6").Value = 23

# Boolean column uses a custom TRUE/FALSE number format
$ws1.Range("C6").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws1.Range("C6").Value = $true

$ws1.Range("D6").Value = "Beer"

# Birthday is a real date formula formatted as mm/dd/yy
$ws1.Range("E6").NumberFormat = "mm/dd/yy"
$ws1.Range("E6").Formula = "=DATE(1940,4,28)"

$ws1.Range("F6").Value = 1.25

# Make the "customers" sheet the active tab/sheet with a new selection
$ws1.Activate() | Out-Null
$ws1.Range("E7").Select() | Out-Null
